$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting existing rows 61..152 down to 62..153
$ws.Rows.Item(61).EntireRow.Insert()

# Populate the new row 61 with data
$ws.Range("A61").Value = 1
$ws.Range("B61").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C61").Value = "Arica y Parinacota"
$ws.Range("D61").Value = 44902
$ws.Range("E61").Value = 15
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100106
$ws.Range("H61").Value = "Oleaginosos"
$ws.Range("I61").Value = 100106002
$ws.Range("J61").Value = "Palta"
$ws.Range("K61").Value = "Edranol"
$ws.Range("L61").Value = "Segunda"
$ws.Range("M61").Value = 300
$ws.Range("N61").Value = 64000
$ws.Range("O61").Value = 65000
$ws.Range("P61").Value = 64500
$ws.Range("Q61").Value = "`$/caja 25 kilos"
$ws.Range("R61").Value = "Región de Coquimbo"
$ws.Range("S61").Value = 2580
$ws.Range("T61").Value = 25

# Ensure the date cell keeps the date number format (style 2 in original workbook)
$ws.Range("D61").NumberFormat = "YYYY-MM-DD HH:MM:SS"
